# Update the test_set_v2 worksheet with the new model v2 results.
# This rewrites A2:B60 (SAMPLE_ID / labels_v1 columns) with the updated
# sample order/labels, and appends a new row (60) that did not exist before.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Sample_0199", 0),
    @("Sample_0172", 0),
    @("Sample_0024", 0),
    @("Sample_0133", 0),
    @("Sample_0214", 1),
    @("Sample_0245", 1),
    @("Sample_0123", 0),
    @("Sample_0106", 1),
    @("Sample_0026", 1),
    @("Sample_0229", 1),
    @("Sample_0216", 1),
    @("Sample_0023", 0),
    @("Sample_0221", 1),
    @("Sample_0111", 0),
    @("Sample_0174", 1),
    @("Sample_0173", 0),
    @("Sample_0205", 0),
    @("Sample_0016", 1),
    @("Sample_0128", 1),
    @("Sample_0247", 1),
    @("Sample_0171", 0),
    @("Sample_0151", 0),
    @("Sample_0198", 1),
    @("Sample_0092", 1),
    @("Sample_0143", 0),
    @("Sample_0120", 0),
    @("Sample_0188", 1),
    @("Sample_0233", 1),
    @("Sample_0045", 1),
    @("Sample_0122", 1),
    @("Sample_0140", 0),
    @("Sample_0179", 1),
    @("Sample_0243", 1),
    @("Sample_0035", 1),
    @("Sample_0197", 1),
    @("Sample_0153", 0),
    @("Sample_0107", 0),
    @("Sample_0027", 0),
    @("Sample_0182", 1),
    @("Sample_0134", 0),
    @("Sample_0012", 0),
    @("Sample_0118", 1),
    @("Sample_0196", 1),
    @("Sample_0155", 0),
    @("Sample_0228", 1),
    @("Sample_0183", 1),
    @("Sample_0225", 1),
    @("Sample_0019", 1),
    @("Sample_0062", 1),
    @("Sample_0222", 1),
    @("Sample_0089", 0),
    @("Sample_0119", 1),
    @("Sample_0195", 1),
    @("Sample_0136", 0),
    @("Sample_0178", 0),
    @("Sample_0184", 0),
    @("Sample_0044", 1),
    @("Sample_0192", 1),
    @("Sample_0230", 1)
)

$styleSource = $ws.Range("A2")

for ($i = 0; $i -lt $data.Length; $i++) {
    $rowNum = $i + 2
    $sampleId = $data[$i][0]
    $label = $data[$i][1]

    $ws.Cells.Item($rowNum, 1).Value = $sampleId
    $ws.Cells.Item($rowNum, 2).Value = $label
}

# Row 60 is a brand new row; make sure column A carries the same style
# (bold/border/centered) as the rest of the SAMPLE_ID column.
$styleSource.Copy()
$ws.Range("A60").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
